$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.183329105377197
$ws.Range("B1").Value = 5.314236640930176
$ws.Range("C1").Value = 2.133415222167969
$ws.Range("D1").Value = 1.294260382652283
$ws.Range("E1").Value = 1.258863687515259
